$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the first four data rows with the new control-point values
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = 62

$ws.Range("A3").Value = 21
$ws.Range("B3").Value = 52

$ws.Range("A4").Value = 12
$ws.Range("B4").Value = 44

$ws.Range("A5").Value = 22
$ws.Range("B5").Value = 6

# Remove the trailing rows (6-9) that are no longer part of the data set
$ws.Range("A6:B9").EntireRow.Delete()
